$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# row 19
$ws.Range("H19").Value = 1160.3478
$ws.Range("I19").Value = 1230.9231
$ws.Range("K19").Value = 1230.9231
$ws.Range("M19").Value = -1055.9231
# row 92
$ws.Range("H92").Value = 447.875
$ws.Range("I92").Value = 511.33334
$ws.Range("K92").Value = 511.33334
$ws.Range("M92").Value = 736.66666
# row 98
$ws.Range("H98").Value = 761
$ws.Range("I98").Value = 805.55554
$ws.Range("J98").Value = 660.75
$ws.Range("K98").Value = 805.55554
$ws.Range("L98").Value = 660.75
$ws.Range("M98").Value = 692.44446
$ws.Range("N98").Value = -3656.75
# row 122
$ws.Range("H122").Value = 761
$ws.Range("I122").Value = 805.55554
$ws.Range("J122").Value = 660.75
$ws.Range("K122").Value = 2416.66662
$ws.Range("L122").Value = 1982.25
$ws.Range("M122").Value = 33.33338000000003
$ws.Range("N122").Value = -6882.25
# row 129
$ws.Range("H129").Value = 2606.625
$ws.Range("J129").Value = 3923.3333
$ws.Range("L129").Value = 11769.9999
$ws.Range("N129").Value = -21769.9999
# row 131
$ws.Range("H131").Value = 2034.375
$ws.Range("I131").Value = 2034.375
$ws.Range("K131").Value = 6103.125
$ws.Range("M131").Value = -1063.125

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# row 32
$ws.Range("H32").Value = 14134.32
$ws.Range("I32").Value = 14698.228
$ws.Range("K32").Value = 14698.228
$ws.Range("M32").Value = -14411.228
# row 96
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()
# row 132
$ws.Range("H132").Value = 3469.342
$ws.Range("I132").Value = 2530.2812
$ws.Range("K132").Value = 7590.8436
$ws.Range("M132").Value = -5060.8436
# row 139
$ws.Range("H139").Value = 73000
$ws.Range("J139").Value = 78750
$ws.Range("L139").Value = 78750
$ws.Range("N139").Value = -89030

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# row 6
$ws.Range("H6").Value = 24995.125
$ws.Range("J6").Value = 24995.125
$ws.Range("L6").Value = 24995.125
$ws.Range("N6").Value = -25221.125
# row 20
$ws.Range("H20").Value = 998
$ws.Range("I20").Value = 998
$ws.Range("K20").Value = 998
$ws.Range("M20").Value = -751
# row 22
$ws.Range("H22").Value = 850.55554
$ws.Range("J22").Value = 715
$ws.Range("L22").Value = 715
$ws.Range("N22").Value = -1061
# row 46
$ws.Range("H46").Value = 13353
$ws.Range("I46").Value = 12529.5
$ws.Range("J46").Value = 15000
$ws.Range("K46").Value = 12529.5
$ws.Range("L46").Value = 15000
$ws.Range("M46").Value = -12231.5
$ws.Range("N46").Value = -15596
# row 80
$ws.Range("H80").Value = 3544.1667
$ws.Range("I80").Value = 346.9091
$ws.Range("J80").Value = 8568.429
$ws.Range("K80").Value = 346.9091
$ws.Range("L80").Value = 8568.429
$ws.Range("M80").Value = 651.0908999999999
$ws.Range("N80").Value = -10564.429
# row 83
$ws.Range("H83").Value = 3544.1667
$ws.Range("I83").Value = 346.9091
$ws.Range("J83").Value = 8568.429
$ws.Range("K83").Value = 1734.5455
$ws.Range("L83").Value = 42842.145
$ws.Range("M83").Value = 3257.4545
$ws.Range("N83").Value = -52826.145
# row 107
$ws.Range("H107").Value = 1206.5385
$ws.Range("I107").Value = 1100.909
$ws.Range("K107").Value = 1100.909
$ws.Range("M107").Value = 819.0909999999999
# row 135
$ws.Range("H135").Value = 99995.164
$ws.Range("J135").Value = 99995.164
$ws.Range("L135").Value = 99995.164
$ws.Range("N135").Value = -110135.164

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# row 4
$ws.Range("H4").Value = 46.555557
$ws.Range("I4").Value = 46.555557
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 46.555557
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 65.44444300000001
$ws.Range("N4").ClearContents()
# row 31
$ws.Range("H31").Value = 8001.1665
$ws.Range("J31").Value = 4555.294
$ws.Range("L31").Value = 4555.294
$ws.Range("N31").Value = -5145.294
# row 34
$ws.Range("H34").Value = 8001.1665
$ws.Range("J34").Value = 4555.294
$ws.Range("L34").Value = 4555.294
$ws.Range("N34").Value = -4959.294
# row 99
$ws.Range("H99").Value = 3085.6
$ws.Range("I99").Value = 3182
$ws.Range("K99").Value = 3182
$ws.Range("M99").Value = -1684
# row 126
$ws.Range("H126").Value = 3085.6
$ws.Range("I126").Value = 3182
$ws.Range("K126").Value = 9546
$ws.Range("M126").Value = -7076

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# row 15
$ws.Range("H15").Value = 339.4
$ws.Range("J15").Value = 174.25
$ws.Range("L15").Value = 522.75
$ws.Range("N15").Value = -802.75
# row 26
$ws.Range("H26").Value = 1500
$ws.Range("J26").Value = 1500
$ws.Range("L26").Value = 4500
$ws.Range("N26").Value = -5076
# row 94
$ws.Range("H94").Value = 8784.143
$ws.Range("I94").Value = 496.33334
$ws.Range("J94").Value = 15000
$ws.Range("K94").Value = 1489.00002
$ws.Range("L94").Value = 45000
$ws.Range("M94").Value = -813.0000199999999
$ws.Range("N94").Value = -46352
# row 100
$ws.Range("H100").Value = 209.5
$ws.Range("I100").Value = 209.5
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 628.5
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = 182.5
$ws.Range("N100").ClearContents()
# row 107
$ws.Range("H107").Value = 381.1613
$ws.Range("I107").Value = 358.45456
$ws.Range("J107").Value = 393.65
$ws.Range("K107").Value = 1075.36368
$ws.Range("L107").Value = 1180.95
$ws.Range("M107").Value = 844.6363200000001
$ws.Range("N107").Value = -5020.95
# row 121
$ws.Range("H121").Value = 3085.1667
$ws.Range("J121").Value = 3562.2
$ws.Range("L121").Value = 10686.6
$ws.Range("N121").Value = -13306.6
# row 131
$ws.Range("H131").Value = 2324.5806
$ws.Range("J131").Value = 2268.7334
$ws.Range("L131").Value = 6806.2002
$ws.Range("N131").Value = -16886.2002

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# row 70
$ws.Range("H70").Value = 2050
$ws.Range("I70").Value = 2050
$ws.Range("K70").Value = 2050
$ws.Range("M70").Value = -1780
# row 73
$ws.Range("H73").Value = 2050
$ws.Range("I73").Value = 2050
$ws.Range("K73").Value = 2050
$ws.Range("M73").Value = -1114
# row 113
$ws.Range("H113").Value = 3072
$ws.Range("I113").Value = 3567.2856
$ws.Range("J113").Value = 1916.3334
$ws.Range("K113").Value = 3567.2856
$ws.Range("L113").Value = 1916.3334
$ws.Range("M113").Value = -1397.2856
$ws.Range("N113").Value = -6256.3334
# row 122
$ws.Range("H122").Value = 433.33334

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# row 2
$ws.Range("H2").Value = 102
$ws.Range("I2").Value = 102
$ws.Range("K2").Value = 102
$ws.Range("M2").Value = 10
# row 16
$ws.Range("H16").Value = 20000
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 20000
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 20000
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -20340
# row 35
$ws.Range("H35").Value = 2631.375
$ws.Range("I35").Value = 2631.375
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 2631.375
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -2295.375
$ws.Range("N35").ClearContents()
# row 93
$ws.Range("H93").Value = 1166.3334
$ws.Range("J93").Value = 1300
$ws.Range("L93").Value = 1300
$ws.Range("N93").Value = -3796
# row 108
$ws.Range("H108").Value = 86649.5
$ws.Range("J108").Value = 86649.5
$ws.Range("L108").Value = 86649.5
$ws.Range("N108").Value = -94329.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# row 96
$ws.Range("H96").Value = 2366.8572
$ws.Range("I96").Value = 2594.6667
$ws.Range("K96").Value = 2594.6667
$ws.Range("M96").Value = -1221.6667
# row 101
$ws.Range("H101").Value = 21993.8
$ws.Range("J101").Value = 21993.8
$ws.Range("L101").Value = 21993.8
$ws.Range("N101").Value = -28483.8
# row 132
$ws.Range("H132").Value = 100001770
$ws.Range("I132").Value = 1784.5714
$ws.Range("K132").Value = 5353.7142
$ws.Range("M132").Value = -2823.7142
